# Add a formula to spreadsheet
# Sheet1!B2 previously held the static value 654; replace it with a formula
# that doubles Sheet1!B1 (456 * 2 = 912), and move the active selection to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Formula = "=B1*2"

$ws.Range("B1").Select() | Out-Null
